$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B header changes from "Grupo" to "Perfil"
$ws.Range("B1").Value = "Perfil"

# Column B values: previously numeric group codes (194/195/196) paired with
# a numeric profile code in column C (3 = Usuario, 4 = Supervisor).
# Now column B directly holds the profile name text, and column C is dropped.
$ws.Range("B2:B7").Value = "Usuario"
$ws.Range("B8:B11").Value = "Supervisor"

# Drop the old column C (Perfil numeric code) entirely - data now lives in B.
$ws.Range("C:C").Delete()

# Restore the selection to the top of the sheet.
$ws.Range("B1").Select() | Out-Null
